# Insert a new data row above row 546 (shifting existing rows 546:619 down
# to 547:620) and populate the new row with the Feria Lagunitas de Puerto
# Montt / Zanahoria record for the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 546 downward by inserting a new blank row at 546.
$ws.Rows.Item(546).Insert()

# Populate the newly inserted row 546 with the new observation's data.
$ws.Range("A546").Value = 4
$ws.Range("B546").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C546").Value = "Los Lagos"
$ws.Range("D546").Value = 45154
$ws.Range("E546").Value = 10
$ws.Range("F546").Value = 100114013
$ws.Range("G546").Value = "Zanahoria"
$ws.Range("H546").Value = "Sin especificar"
$ws.Range("I546").Value = "Primera"
$ws.Range("J546").Value = 150
$ws.Range("K546").Value = 7500
$ws.Range("L546").Value = 7500
$ws.Range("M546").Value = 7500
$ws.Range("N546").Value = "$/saco 20 kilos"
$ws.Range("O546").Value = "Provincia de Llanquihue"
$ws.Range("P546").Value = 375
$ws.Range("Q546").Value = 20
$ws.Range("R546").Value = "Hortaliza"
